$d = $word.ActiveDocument

function Get-TargetRange {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*allerede viste*") {
            return $p.Range
        }
    }
    return $null
}

# Force a clean run-split (without altering the visible formatting) on the
# range [start, end) by toggling Bold on then back off. Word always breaks
# the containing run at both edges of the formatted sub-range, giving us a
# dedicated <w:r> for exactly that text, with the original formatting intact.
function Split-Range {
    param([int]$start, [int]$end)
    $rng = $d.Range($start, $end)
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------------
# Phase 1: fix up the sentence's wording
# ---------------------------------------------------------------------------

# Remove the now-redundant " er derfor valgt ikke at udarbejde" that currently
# follows "charts" - do this first (before the phrase is reinserted below)
# so the Find cannot match the wrong occurrence.
$rng1 = Get-TargetRange
$rng1.Find.Execute(
    " er derfor valgt ikke at udarbejde",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2
)

# "... viste. Flow " -> "... viste. Der er derfor valgt ikke at udarbejde flow "
$rng2 = Get-TargetRange
$rng2.Find.Execute(
    "ikke allerede viste. Flow ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ikke allerede viste. Der er derfor valgt ikke at udarbejde flow ",
    2
)

# ---------------------------------------------------------------------------
# Phase 2: reproduce the run layout
# ---------------------------------------------------------------------------

$full = Get-TargetRange
$base = $full.Start

# "ikke allerede viste." ends here (length of "Process view ... ikke allerede viste.")
$cursor = $d.Range($base, $full.End)
$cursor.Find.Execute("ikke allerede viste.", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$posEndViste = $cursor.End

# " Der er derfor valgt ikke at udarbejde flow " follows; find "udarbejde flow "
$cursor2 = $d.Range($posEndViste, $full.End)
$cursor2.Find.Execute("udarbejde flow ", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$posUdarbejdeStart = $cursor2.Start
$posUdarbejdeEnd = $posUdarbejdeStart + ("udarbejde").Length   # end of "udarbejde"
$posSpaceEnd = $posUdarbejdeEnd + 1                             # end of the single space before "flow"
$posFEnd = $posSpaceEnd + 1                                     # end of "f"
$posLowEnd = $cursor2.End                                       # end of "low "

# Split 1: " ikke allerede viste." | " Der er derfor valgt ikke at udarbejde"
Split-Range $base $posEndViste
Split-Range $posEndViste $posUdarbejdeEnd

# Split 2: " " | "f" | "low "
Split-Range $posUdarbejdeEnd $posSpaceEnd
Split-Range $posSpaceEnd $posFEnd
Split-Range $posFEnd $posLowEnd

# Now locate the comma + " da sekvensdiagrammerne" + " lettere..." region.
$cursor3 = $d.Range($posLowEnd, $full.End)
$cursor3.Find.Execute(", da sekvensdiagrammerne lettere", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$posCommaStart = $cursor3.Start
$posCommaEnd = $posCommaStart + 1                                # end of ","
$posSekvensEnd = $posCommaEnd + (" da sekvensdiagrammerne").Length

Split-Range $posCommaStart $posCommaEnd
Split-Range $posCommaEnd $posSekvensEnd

# ---------------------------------------------------------------------------
# Phase 3: relocate the _GoBack bookmark from the end of the paragraph to
# right after " da sekvensdiagrammerne" (before " lettere kan relateres...").
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($posSekvensEnd, $posSekvensEnd))

Write-Output (Get-TargetRange).Text
